$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "SamplesTab" query (B3) to use samp.sample_tumor_status directly
# instead of the collected "tumor" alias, and tweak the ORDER BY indentation.
$newTumorQuery = @'
MATCH (s:study)<--(p:participant)<--(samp:sample)
WHERE s.study_name in ["CIDR: Discovery, Biology, and Risk of Inherited Variants in Glioma sample"]
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as `Sample ID`,
 coalesce(p.participant_id,'') as `Participant ID`,
 coalesce(s.study_name, '') as `Study Name`,
 coalesce(s.phs_accession,'') as `Accession`,
 coalesce(samp.sample_tumor_status,'') as `Tumor`,
coalesce(samp.sample_type,'') as `Analyte Type`
  ORDER By samp.sample_id LIMIT 100
'@

$ws.Range("B3").Value = $newTumorQuery

# Move the active selection from E4 to D11
$ws.Range("D11").Select()
